$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '71.152.83'
Set-TextValue 'E2' '  -0.52%  '
Set-TextValue 'D3' '3.833.29'
Set-TextValue 'E3' '  +0.22%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '704.30'
Set-TextValue 'E5' '  +0.42%  '
Set-TextValue 'D6' '172.07'
Set-TextValue 'E6' '  -1.20%  '
Set-TextValue 'D7' '3.831.83'
Set-TextValue 'E7' '  +0.26%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'E9' '  -0.56%  '
Set-TextValue 'E10' '  -0.54%  '
Set-TextValue 'D11' '7.38'
Set-TextValue 'E11' '  -0.64%  '
Set-TextValue 'E12' '  -0.48%  '
Set-TextValue 'E13' '  -1.74%  '
Set-TextValue 'D14' '36.73'
Set-TextValue 'E14' '  +0.35%  '
Set-TextValue 'D15' '4.481.94'
Set-TextValue 'E15' '  +0.48%  '
Set-TextValue 'D16' '3.808.46'
Set-TextValue 'E16' '  +0.43%  '
Set-TextValue 'D17' '71.136.64'
Set-TextValue 'E17' '  -0.45%  '
Set-TextValue 'E18' '  -0.08%  '
Set-TextValue 'E19' '  +0.28%  '
Set-TextValue 'E20' '  -2.20%  '
Set-TextValue 'E21' '  -3.52%  '
Set-TextValue 'D22' '495.38'
Set-TextValue 'E22' '  +1.88%  '
Set-TextValue 'D23' '0.738'
Set-TextValue 'E23' '  +2.69%  '
Set-TextValue 'D24' '85.38'
Set-TextValue 'E24' '  +0.69%  '
Set-TextValue 'E25' '  +0.54%  '
Set-TextValue 'D26' '10.64'
Set-TextValue 'E26' '  +0.81%  '
Set-TextValue 'E27' '  -2.11%  '
Set-TextValue 'E28' '  -2.90%  '
Set-TextValue 'E29' '  -0.06%  '
Set-TextValue 'D30' '3.10'
Set-TextValue 'E30' '  -2.35%  '
Set-TextValue 'D31' '7.45'
Set-TextValue 'E31' '  -2.18%  '
Set-TextValue 'E32' '  -3.53%  '
Set-TextValue 'D33' '29.45'
Set-TextValue 'E33' '  -1.03%  '
Set-TextValue 'E34' '  -3.63%  '
Set-TextValue 'E35' '  -1.12%  '
Set-TextValue 'B36' 'Binance-PegBSC-USD'
Set-TextValue 'C36' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D36' '1.01'
Set-TextValue 'E36' '  +1.03%  '
Set-TextValue 'B37' 'RenzoRestakedETH'
Set-TextValue 'C37' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D37' '3.796.64'
Set-TextValue 'E37' '  +0.68%  '
Set-TextValue 'E38' '  -0.92%  '
Set-TextValue 'D39' '2.35'
Set-TextValue 'E39' '  -2.27%  '
Set-TextValue 'D40' '1.04'
Set-TextValue 'E40' '  +4.44%  '
Set-TextValue 'E41' '  -0.79%  '
Set-TextValue 'E42' '  -3.09%  '
Set-TextValue 'E43' '  -0.01%  '
Set-TextValue 'E44' '  +0.18%  '
Set-TextValue 'E45' '  +0.84%  '
Set-TextValue 'D46' '163.67'
Set-TextValue 'E46' '  +0.21%  '
Set-TextValue 'D47' '429.08'
Set-TextValue 'E47' '  +3.66%  '
Set-TextValue 'D48' '48.91'
Set-TextValue 'E48' '  +0.35%  '
Set-TextValue 'D49' '8.78'
Set-TextValue 'E49' '  +0.83%  '
Set-TextValue 'E50' '  -0.23%  '
Set-TextValue 'E51' '  -1.90%  '
